$d = $word.ActiveDocument

# --- Update the "FECHA DE ENTREGA" title paragraph date ---
$d.Content.Find.Execute("15/10/2025", $true, $false, $false, $false, $false, $true, 1, $false, "21/10/2025", 2) | Out-Null

# --- Update the Sprint Backlog table dates ---
$t = $d.Tables.Item(1)

$t.Cell(3, 4).Range.Find.Execute("21/10/2025", $true, $false, $false, $false, $false, $true, 1, $false, "14/10/2025", 2) | Out-Null
$t.Cell(3, 5).Range.Find.Execute("24/10/2025", $true, $false, $false, $false, $false, $true, 1, $false, "17/10/2025", 2) | Out-Null
$t.Cell(4, 4).Range.Find.Execute("25/10/2025", $true, $false, $false, $false, $false, $true, 1, $false, "18/10/2025", 2) | Out-Null
$t.Cell(4, 5).Range.Find.Execute("28/10/2025", $true, $false, $false, $false, $false, $true, 1, $false, "21/10/2025", 2) | Out-Null
$t.Cell(5, 4).Range.Find.Execute("29/10/2025", $true, $false, $false, $false, $false, $true, 1, $false, "21/10/2025", 2) | Out-Null
$t.Cell(5, 5).Range.Find.Execute("04/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "26/10/2025", 2) | Out-Null
$t.Cell(6, 4).Range.Find.Execute("05/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "27/10/2025", 2) | Out-Null
$t.Cell(6, 5).Range.Find.Execute("05/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "28/10/2025", 2) | Out-Null
$t.Cell(7, 4).Range.Find.Execute("06/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "28/10/2025", 2) | Out-Null
$t.Cell(7, 5).Range.Find.Execute("11/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "02/11/2025", 2) | Out-Null
$t.Cell(8, 4).Range.Find.Execute("12/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "03/11/2025", 2) | Out-Null
$t.Cell(8, 5).Range.Find.Execute("13/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "04/11/2025", 2) | Out-Null
$t.Cell(9, 4).Range.Find.Execute("14/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "04/11/2025", 2) | Out-Null
$t.Cell(9, 5).Range.Find.Execute("17/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "08/11/2025", 2) | Out-Null
$t.Cell(10, 4).Range.Find.Execute("18/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "09/11/2025", 2) | Out-Null
$t.Cell(10, 5).Range.Find.Execute("21/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "11/11/2025", 2) | Out-Null
$t.Cell(11, 4).Range.Find.Execute("22/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "11/11/2025", 2) | Out-Null
$t.Cell(11, 5).Range.Find.Execute("26/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "15/11/2025", 2) | Out-Null
$t.Cell(12, 4).Range.Find.Execute("27/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "16/11/2025", 2) | Out-Null
$t.Cell(12, 5).Range.Find.Execute("29/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "18/11/2025", 2) | Out-Null
$t.Cell(13, 4).Range.Find.Execute("30/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "18/11/2025", 2) | Out-Null
$t.Cell(13, 5).Range.Find.Execute("03/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "22/11/2025", 2) | Out-Null
$t.Cell(14, 4).Range.Find.Execute("04/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "23/11/2025", 2) | Out-Null
$t.Cell(14, 5).Range.Find.Execute("07/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "25/11/2025", 2) | Out-Null
$t.Cell(15, 4).Range.Find.Execute("08/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "25/11/2025", 2) | Out-Null
$t.Cell(15, 5).Range.Find.Execute("13/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "29/12/2025", 2) | Out-Null
$t.Cell(16, 4).Range.Find.Execute("14/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "30/12/2025", 2) | Out-Null
$t.Cell(16, 5).Range.Find.Execute("15/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "2/12/2025", 2) | Out-Null
$t.Cell(17, 4).Range.Find.Execute("16/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "2/12/2025", 2) | Out-Null
$t.Cell(17, 5).Range.Find.Execute("18/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "04/12/2025", 2) | Out-Null
$t.Cell(18, 4).Range.Find.Execute("19/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "05/12/2025", 2) | Out-Null
$t.Cell(18, 5).Range.Find.Execute("21/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "07/12/2025", 2) | Out-Null
$t.Cell(19, 4).Range.Find.Execute("22/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "08/12/2025", 2) | Out-Null
$t.Cell(19, 5).Range.Find.Execute("23/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "09/12/2025", 2) | Out-Null
$t.Cell(20, 4).Range.Find.Execute("24/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "09/12/2025", 2) | Out-Null
$t.Cell(20, 5).Range.Find.Execute("28/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "11/12/2025", 2) | Out-Null
$t.Cell(21, 4).Range.Find.Execute("29/11/2025", $true, $false, $false, $false, $false, $true, 1, $false, "12/11/2025", 2) | Out-Null
$t.Cell(21, 5).Range.Find.Execute("30/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "14/12/2025", 2) | Out-Null
$t.Cell(22, 4).Range.Find.Execute("31/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "15/12/2025", 2) | Out-Null
$t.Cell(22, 5).Range.Find.Execute("31/12/2025", $true, $false, $false, $false, $false, $true, 1, $false, "16/12/2025", 2) | Out-Null
$t.Cell(23, 4).Range.Find.Execute("01/01/2026", $true, $false, $false, $false, $false, $true, 1, $false, "16/12/2026", 2) | Out-Null
$t.Cell(23, 5).Range.Find.Execute("07/01/2026", $true, $false, $false, $false, $false, $true, 1, $false, "23/12/2026", 2) | Out-Null
$t.Cell(24, 4).Range.Find.Execute("08/01/2026", $true, $false, $false, $false, $false, $true, 1, $false, "23/12/2026", 2) | Out-Null
$t.Cell(24, 5).Range.Find.Execute("13/01/2026", $true, $false, $false, $false, $false, $true, 1, $false, "27/12/2026", 2) | Out-Null
$t.Cell(25, 4).Range.Find.Execute("14/01/2026", $true, $false, $false, $false, $false, $true, 1, $false, "28/12/2026", 2) | Out-Null
$t.Cell(25, 5).Range.Find.Execute("15/01/2026", $true, $false, $false, $false, $false, $true, 1, $false, "30/12/2026", 2) | Out-Null
